$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data pull.
# D-column values are entered with a leading apostrophe so Excel stores them
# as text (matching the source data's formatting, e.g. trailing zeros like '1.00').

$ws.Range('D2').Formula = "'51.498.70"
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Formula = "'3.018.91"
$ws.Range('E3').Value = '  +3.36%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Formula = "'379.39"
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').Formula = "'102.53"
$ws.Range('E6').Value = '  +3.03%  '
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('D9').Formula = "'0.591"
$ws.Range('E9').Value = '  +3.57%  '
$ws.Range('D10').Formula = "'36.61"
$ws.Range('E10').Value = '  +3.07%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Formula = "'0.0855"
$ws.Range('D13').Formula = "'3.498.22"
$ws.Range('E13').Value = '  +3.41%  '
$ws.Range('D14').Formula = "'18.45"
$ws.Range('E14').Value = '  +2.80%  '
$ws.Range('D15').Formula = "'7.72"
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Formula = "'3.017.90"
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').Formula = "'0.982"
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Formula = "'10.31"
$ws.Range('E18').Value = '  -13.88%  '
$ws.Range('D19').Formula = "'51.541.81"
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').Formula = "'0.0₃0960"
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').Formula = "'70.03"
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('D24').Formula = "'267.27"
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').Formula = "'3.15"
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').Formula = "'8.18"
$ws.Range('E26').Value = '  +4.57%  '
$ws.Range('D27').Formula = "'7.49"
$ws.Range('E27').Value = '  +6.42%  '
$ws.Range('D28').Formula = "'0.171"
$ws.Range('E28').Value = '  +6.35%  '
$ws.Range('D29').Formula = "'1.00"
$ws.Range('D30').Formula = "'26.14"
$ws.Range('E30').Value = '  +3.29%  '
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('D32').Formula = "'10.28"
$ws.Range('E32').Value = '  +3.55%  '
$ws.Range('D33').Formula = "'34.11"
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('D34').Formula = "'50.56"
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').Formula = "'2.05"
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('E36').Value = '  +5.11%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').Formula = "'3.26"
$ws.Range('E38').Value = '  +6.90%  '
$ws.Range('D39').Formula = "'17.29"
$ws.Range('E39').Value = '  +6.51%  '
$ws.Range('E40').Value = '  +4.32%  '
$ws.Range('D41').Formula = "'2.59"
$ws.Range('E41').Value = '  +7.79%  '
$ws.Range('D42').Formula = "'0.281"
$ws.Range('E42').Value = '  +10.09%  '
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Formula = "'126.51"
$ws.Range('E44').Value = '  +2.81%  '
$ws.Range('D45').Formula = "'3.72"
$ws.Range('E45').Value = '  +10.83%  '
$ws.Range('D46').Formula = "'21.98"
$ws.Range('E46').Value = '  +5.41%  '
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').Formula = "'2.029.68"
$ws.Range('D50').Formula = "'3.320.28"
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('E51').Value = '  +1.93%  '
